$d = $word.ActiveDocument

# 1. "People who travel" paragraph: merge the three runs (incl. the "on a daily basis" run
#    wrapped in gramStart/gramEnd proofErr markers) into a single run with no proofErr.
$d.Content.Find.Execute(
    "People who travel on a daily basis or even if someone new commuting to the highway will have an idea of the road condition and how time, weather and other scenarios while driving may cause collision. They will get awareness to take precautions",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "People who travel on a daily basis or even if someone new commuting to the highway will have an idea of the road condition and how time, weather and other scenarios while driving may cause collision. They will get awareness to take precautions",
    2)

# 2. "Numpy" list item: drop the spellStart/spellEnd proofErr wrapper around it. Spanning the
#    replace across the preceding paragraph mark forces Word to regenerate the run cleanly.
$d.Content.Find.Execute("Pandas" + [char]13 + "Numpy", $true, $false, $false, $false, $false, $true, 1, $false, "Pandas" + [char]13 + "Numpy", 2)

# 3. "Folium" -> "Matplotlib"
$d.Content.Find.Execute("Folium", $true, $false, $false, $false, $false, $true, 1, $false, "Matplotlib", 2)

# 4. Remove the "Geocode" list item paragraph entirely.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Geocode") {
        $p.Range.Delete()
    }
}

# 5. Remove the now-duplicate "Matplotlib" list item paragraph (the one that used to follow
#    "Geocode"); the first "Matplotlib" (ex-"Folium") stays.
$matplotlibParas = @()
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13) -eq "Matplotlib") {
        $matplotlibParas += $p
    }
}
if ($matplotlibParas.Count -gt 1) {
    $matplotlibParas[1].Range.Delete()
}

# 6. "Scikit" + "-learn": merge into a single run "Scikit-learn" with no proofErr wrapper.
$d.Content.Find.Execute("Matplotlib" + [char]13 + "Scikit-learn", $true, $false, $false, $false, $false, $true, 1, $false, "Matplotlib" + [char]13 + "Scikit-learn", 2)

# 7. "As a database..." paragraph: merge all runs into one, dropping the proofErr wrappers
#    around "Collission" and "on the basis of" (the visible text itself is unchanged).
$dbText = "As a database, I used GitHub repository in my study. My master data which has the main components Severity Code, Weather Elements, Light Condition, Junction Type and Collission Type. For classification problems, matplotlib is extensively used the supervised algorithm. The primary perspective of this algorithm is predicting the value of the desired variable by learning decision rules deduced from the features of the data and create a model of that. A root node is designated for the construction of this model based on the best attribute picked by the gain approach and the sub-nodes are then generated on the basis of the decision taken in relation to the status of quality selected at each node. When each node is reduced to a single quality status, the class is determined at the end of the node; it is called a leaf. These courses of action continue recursively until a class is defined at the end of each node."
$d.Content.Find.Execute($dbText, $true, $false, $false, $false, $false, $true, 1, $false, $dbText, 2)

# 8. Results paragraph: "Ada-Boost" -> "Gradient-Boost and Random Forest", split across
#    separate runs the same way Word would after a few distinct typing/editing passes, and
#    move the "_GoBack" last-edit-position bookmark to sit right before "gives the best result".
$rng = $d.Content
$rng.Find.Execute("Ada-Boost ") | Out-Null
$start = $rng.Start
$rng.Text = "Gradient-Boost and Random Forest "

$splitOffsets = @(0, 8, 15)   # before "Gradient", before "-Boost ", before "and Random Forest "
foreach ($off in $splitOffsets) {
    $pos = $start + $off
    $bm = $d.Range($pos, $pos)
    $d.Bookmarks.Add("tempmark", $bm)
    $d.Bookmarks.Item("tempmark").Delete()
}

$bmPos = $start + 33   # right before "gives the best result..."
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 9. Discussion paragraph: merge runs into one, dropping the gramStart/gramEnd proofErr
#    wrapper around "get" (the visible text itself is unchanged).
$discussionText = "We observe that most of the accidents in our dataset are Fatal and value for the other three classes is very low. For that reason, in second experiment, we merge Grievous, Simple Injury, Motor Collision these three accident severity classes into one class. Therefore, we have attained the performances of the proposed approaches for two accident severity classes (Fatal / Grievous). In this experiment, we have noticed that the accuracy of interpolated data get increased and remain the same. But it is also mentionable that, the performance is much better than the previous experiment as precision and F1 score increased here in a noticeable way. Besides this, we did experiment with the features in our dataset and have tried to find out their effect on a traffic accident. Statistically I have found that based on the condition of some features the number of accidents gets increased. It’s a significant noticeable thing for making proper steps to decrease the number of accidents."
$d.Content.Find.Execute($discussionText, $true, $false, $false, $false, $false, $true, 1, $false, $discussionText, 2)

# 10. Conclusion paragraph: merge runs into one, dropping the spellStart/spellEnd proofErr
#     wrapper around "servere" (the visible text itself is unchanged).
$conclusionText = "As a result, people are turning to big cities to start a business or work. From the above table, we can see that servere car accidents occurs frequently under clear and dry condition at intersection. Besides, speeding is also an important factor leading the accident happen."
$d.Content.Find.Execute($conclusionText, $true, $false, $false, $false, $false, $true, 1, $false, $conclusionText, 2)
